$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Price {
    param($addr, $value)
    $ws.Range($addr).NumberFormat = "@"
    $ws.Range($addr).Value = $value
}

function Set-Text {
    param($addr, $value)
    $ws.Range($addr).Value = $value
}

# Row 2 - Bitcoin
Set-Price "D2" "64.286.95"
Set-Text  "E2" "  +1.42%  "

# Row 3 - Ethereum
Set-Price "D3" "3.096.51"
Set-Text  "E3" "  +1.14%  "

# Row 4 - TetherUSD (price unchanged)
Set-Text  "E4" "  +0.10%  "

# Row 5 - BNB
Set-Price "D5" "560.06"
Set-Text  "E5" "  +1.86%  "

# Row 6 - Solana
Set-Price "D6" "144.64"
Set-Text  "E6" "  +2.72%  "

# Row 7 - USDC (price unchanged)
Set-Text  "E7" "  +0.07%  "

# Row 8 - LidoStakedEther
Set-Price "D8" "3.094.72"
Set-Text  "E8" "  +1.27%  "

# Row 9 - XRP
Set-Price "D9" "0.508"
Set-Text  "E9" "  +1.44%  "

# Row 10 - Dogecoin
Set-Price "D10" "0.155"
Set-Text  "E10" "  +2.03%  "

# Row 11 - Toncoin
Set-Price "D11" "6.17"
Set-Text  "E11" "  -5.44%  "

# Row 12 - Cardano
Set-Price "D12" "0.473"
Set-Text  "E12" "  +3.74%  "

# Row 13 - ShibaInu (price unchanged)
Set-Text  "E13" "  +0.88%  "

# Row 14 - Avalanche
Set-Price "D14" "35.23"
Set-Text  "E14" "  +1.23%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-Price "D15" "3.604.68"
Set-Text  "E15" "  +1.24%  "

# Row 16 - WrappedBTC
Set-Price "D16" "64.365.71"
Set-Text  "E16" "  +1.69%  "

# Row 17 - WrappedEther
Set-Price "D17" "3.097.38"
Set-Text  "E17" "  +1.00%  "

# Row 18 - TRON
Set-Price "D18" "0.110"
Set-Text  "E18" "  +1.22%  "

# Row 19 - Polkadot
Set-Price "D19" "6.76"
Set-Text  "E19" "  -0.18%  "

# Row 20 - BitcoinCash
Set-Price "D20" "482.41"
Set-Text  "E20" "  -0.30%  "

# Row 21 - Chainlink
Set-Price "D21" "14.06"
Set-Text  "E21" "  +1.60%  "

# Row 22 - Polygon
Set-Price "D22" "0.679"
Set-Text  "E22" "  +0.31%  "

# Row 23 - Uniswap
Set-Price "D23" "7.57"
Set-Text  "E23" "  +4.17%  "

# Row 24 - InternetComputer(DFINITY)
Set-Price "D24" "14.14"
Set-Text  "E24" "  +11.35%  "

# Row 25 - Litecoin
Set-Price "D25" "81.57"
Set-Text  "E25" "  +0.67%  "

# Row 26 - Dai (price unchanged)
Set-Text  "E26" "  +0.07%  "

# Row 27 - PancakeSwap (price unchanged)
Set-Text  "E27" "  +1.64%  "

# Row 28 - RenderToken
Set-Price "D28" "8.04"
Set-Text  "E28" "  +2.12%  "

# Row 29 - ImmutableX (price unchanged)
Set-Text  "E29" "  +2.44%  "

# Row 30 - FirstDigitalUSD (price unchanged)
Set-Text  "E30" "  +0.08%  "

# Row 31 - EthereumClassic
Set-Price "D31" "26.39"
Set-Text  "E31" "  +0.78%  "

# Row 32 - Mantle
Set-Price "D32" "1.16"
Set-Text  "E32" "  -0.63%  "

# Row 33 - Stacks (price unchanged)
Set-Text  "E33" "  +1.40%  "

# Row 34 - NEARProtocol (price unchanged)
Set-Text  "E34" "  -0.26%  "

# Row 35 - Filecoin
Set-Price "D35" "6.22"
Set-Text  "E35" "  +3.93%  "

# Row 36 - OKB
Set-Price "D36" "55.61"
Set-Text  "E36" "  +0.35%  "

# Row 37/38 - VeChain and dogwifhat swapped places (dogwifhat now ranks above VeChain)
Set-Text  "B37" "dogwifhat"
Set-Text  "C37" "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-Price "D37" "2.99"
Set-Text  "E37" "  +16.88%  "

Set-Text  "B38" "VeChain"
Set-Text  "C38" "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
Set-Price "D38" "0.0409"
Set-Text  "E38" "  +2.92%  "

# Row 39 - Bittensor
Set-Price "D39" "443.51"
Set-Text  "E39" "  -4.55%  "

# Row 40 - Hedera
Set-Price "D40" "0.0816"
Set-Text  "E40" "  -0.88%  "

# Row 41 - Maker
Set-Price "D41" "2.965.63"
Set-Text  "E41" "  -2.66%  "

# Row 42 - Cosmos
Set-Price "D42" "8.23"
Set-Text  "E42" "  +0.11%  "

# Row 43 - Kaspa (price unchanged)
Set-Text  "E43" "  -3.76%  "

# Row 44 - InjectiveProtocol
Set-Price "D44" "28.22"
Set-Text  "E44" "  +1.55%  "

# Row 45 - TheGraph
Set-Price "D45" "0.261"
Set-Text  "E45" "  +2.53%  "

# Row 46/47 - Fetch.AI and USDe swapped places (USDe now ranks above Fetch.AI)
Set-Text  "B46" "USDe"
Set-Text  "C46" "https://coinranking.com/coin/exbfr2U-0+usde-usde"
Set-Price "D46" "1.00"
Set-Text  "E46" "  +0.02%  "

Set-Text  "B47" "Fetch.AI"
Set-Text  "C47" "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
Set-Price "D47" "2.14"
Set-Text  "E47" "  +5.01%  "

# Row 48 - Stellar (price unchanged)
Set-Text  "E48" "  +1.75%  "

# Row 49/50 - PEPE and Monero swapped places (Monero now ranks above PEPE)
Set-Text  "B49" "Monero"
Set-Text  "C49" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-Price "D49" "118.41"
Set-Text  "E49" "  +1.32%  "

Set-Text  "B50" "PEPE"
Set-Text  "C50" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-Price "D50" "0.0₃0518"
Set-Text  "E50" "  +1.76%  "

# Row 51 - ThetaToken
Set-Price "D51" "2.10"
Set-Text  "E51" "  +1.19%  "
